# date and numeric column formatting adjustment
#
# Adds a new "formatting function" column (C) entry for the WHO,
# Dashboard/epi, STATCAN and modeling extract rows on the "instruct"
# sheet, and updates the sheet's active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- WHO rows -------------------------------------------------------
$ws.Range("C4").Value = "who_cols_formatting"
$ws.Range("C5").Value = "who_cols_formatting"

# --- Dashboard / epi rows -------------------------------------------
$ws.Range("C8").Value  = "dashboard_and_epi_cols_formatting"
$ws.Range("C9").Value  = "dashboard_and_epi_cols_formatting"
$ws.Range("C10").Value = "dashboard_and_epi_cols_formatting"
$ws.Range("C11").Value = "dashboard_and_epi_cols_formatting"

# --- STATCAN row ------------------------------------------------------
$ws.Range("C15").Value = "statcan_cols_formatting"

# --- modeling rows ----------------------------------------------------
$ws.Range("C16").Value = "modelling_cols_formatting"
$ws.Range("C17").Value = "modelling_cols_formatting"
$ws.Range("C18").Value = "modelling_cols_formatting"

# --- Update the visible window / selection ----------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

$ws.Range("C15:C18").Select()
